$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a throwaway placeholder row ("thanos") sitting in row 2 while
# the new weight-slider cell was being wired up. Remove that row now that the
# slider feeds real initial weights - this also drops its now-unused label
# from the shared string table so every instrument name below it is correct
# again (PARSTEI LX Equity, LEF1TREU Index, ...).
$ws.Rows(2).Delete()

# Re-populate Initial Weights (B) / Opt Portfolio (C) / Opt Portfolio with
# View (D) for every instrument with the values recalculated from the new
# slider-driven initial weights (column B still sums to 1).
$ws.Range("B2").Value = 0.3
$ws.Range("C2").Value = 0.000000000000000002927345865710861971820122562348842620849609375
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = 0.2
$ws.Range("C3").Value = 0.02124717618409318
$ws.Range("D3").Value = 0.02124715427314639

$ws.Range("B4").Value = 0.15
$ws.Range("C4").Value = 0.1273382014015918
$ws.Range("D4").Value = 0.1273388940376063

$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = 0.1697704453391681
$ws.Range("D5").Value = 0.1697699223641621

$ws.Range("B6").Value = 0.1
$ws.Range("C6").Value = 0.2390021039874927
$ws.Range("D6").Value = 0.2390017823104368

$ws.Range("B7").Value = 0.05
$ws.Range("C7").Value = 0.09744142026537399
$ws.Range("D7").Value = 0.09744136897287403

$ws.Range("B8").Value = 0.1
$ws.Range("C8").Value = 0.34520065282228
$ws.Range("D8").Value = 0.3452008780417743
